$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "remote pair programming" matrix (rows 14-20)

# Felix Schmidt row (14): Felix/Michael pairing removed, Felix/Jakob minutes increased
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = 480

# Ermal Gashi row (15): Ermal/Michael pairing removed, Ermal/Magdalena and
# Ermal/Florian minutes increased
$ws.Range("D15").ClearContents()
$ws.Range("F15").Value = 360
$ws.Range("H15").Value = 720

# Thomas Pinheiro de Souza / Florian Buchacher minutes reduced
$ws.Range("H19").Value = 360

# Florian's total with Thomas (G20) now mirrors H19 via formula instead of a
# hard-coded literal value
$ws.Range("G20").Formula = "=H19"

# Update the selected cell to match the saved view state
$ws.Range("G20").Select()
